# Add a new metric row ("# New Vaccination Doses") to the Metrics table.
# This mirrors inserting a new row inside the "Metrics" table right after
# the existing "# Vaccinated People" row (i.e. at worksheet row 42),
# which pushes all subsequent metric rows down by one and appends a
# fresh blank row at the very end of the (now one-row-larger) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at row 42; this shifts rows 42-60 down to 43-61
#    and carries along cell formatting (styles) from the split row.
$ws.Rows(42).Insert()

# 2) Grow the "Metrics" table (and its AutoFilter) so it covers the new
#    row, i.e. A1:F60 -> A1:F61.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F61"))

# 3) Populate the newly inserted row 42 with the new metric.
$ws.Range("A42").Value2 = "Vaccinations"
$ws.Range("B42").Value2 = 40
$ws.Range("C42").Value2 = "# New Vaccination Doses"
$ws.Range("D42").Value2 = 410
$ws.Range("E42").Value2 = ""
$ws.Range("F42").Value2 = "X"

# 4) The "Metric - Sort" column (D) is a simple running sequence
#    (row-1)*10 that is tied to the row position, not to the metric it
#    describes, so after the insert it needs to be re-sequenced for every
#    row from the inserted row through the new last row.
for ($r = 43; $r -le 61; $r++) {
    $ws.Cells.Item($r, 4).Value2 = ($r - 1) * 10
}

# 5) Restore the selection/active cell as left by the edit.
$ws.Range("D57").Select() | Out-Null

Write-Host "Added '# New Vaccination Doses' metric row; table now spans A1:F61"
